$d = $word.ActiveDocument
$d.Content.Find.Execute("547÷2=273, 1", $true, $false, $false, $false, $false, $true, 1, $false, "699÷7=99, 6", 2) | Out-Null
$d.Content.Find.Execute("814÷6=135, 4", $true, $false, $false, $false, $false, $true, 1, $false, "582÷3=194, 0", 2) | Out-Null
$d.Content.Find.Execute("982÷9=109, 1", $true, $false, $false, $false, $false, $true, 1, $false, "527÷2=263, 1", 2) | Out-Null
$d.Content.Find.Execute("929÷6=154, 5", $true, $false, $false, $false, $false, $true, 1, $false, "996÷5=199, 1", 2) | Out-Null
$d.Content.Find.Execute("306÷3=102, 0", $true, $false, $false, $false, $false, $true, 1, $false, "121÷8=15, 1", 2) | Out-Null
$d.Content.Find.Execute("202÷8=25, 2", $true, $false, $false, $false, $false, $true, 1, $false, "707÷3=235, 2", 2) | Out-Null
$d.Content.Find.Execute("260÷8=32, 4", $true, $false, $false, $false, $false, $true, 1, $false, "611÷2=305, 1", 2) | Out-Null
$d.Content.Find.Execute("125÷7=17, 6", $true, $false, $false, $false, $false, $true, 1, $false, "990÷9=110, 0", 2) | Out-Null
$d.Content.Find.Execute("981÷6=163, 3", $true, $false, $false, $false, $false, $true, 1, $false, "488÷9=54, 2", 2) | Out-Null
$d.Content.Find.Execute("685÷2=342, 1", $true, $false, $false, $false, $false, $true, 1, $false, "986÷8=123, 2", 2) | Out-Null
$d.Content.Find.Execute("478÷5=95, 3", $true, $false, $false, $false, $false, $true, 1, $false, "308÷5=61, 3", 2) | Out-Null
$d.Content.Find.Execute("105÷5=21, 0", $true, $false, $false, $false, $false, $true, 1, $false, "312÷3=104, 0", 2) | Out-Null
$d.Content.Find.Execute("336÷4=84, 0", $true, $false, $false, $false, $false, $true, 1, $false, "838÷9=93, 1", 2) | Out-Null
$d.Content.Find.Execute("557÷8=69, 5", $true, $false, $false, $false, $false, $true, 1, $false, "902÷5=180, 2", 2) | Out-Null
$d.Content.Find.Execute("914÷7=130, 4", $true, $false, $false, $false, $false, $true, 1, $false, "145÷8=18, 1", 2) | Out-Null
$d.Content.Find.Execute("776÷8=97, 0", $true, $false, $false, $false, $false, $true, 1, $false, "905÷5=181, 0", 2) | Out-Null
$d.Content.Find.Execute("773÷5=154, 3", $true, $false, $false, $false, $false, $true, 1, $false, "628÷2=314, 0", 2) | Out-Null
$d.Content.Find.Execute("953÷9=105, 8", $true, $false, $false, $false, $false, $true, 1, $false, "880÷7=125, 5", 2) | Out-Null
$d.Content.Find.Execute("212÷3=70, 2", $true, $false, $false, $false, $false, $true, 1, $false, "589÷2=294, 1", 2) | Out-Null
$d.Content.Find.Execute("247÷6=41, 1", $true, $false, $false, $false, $false, $true, 1, $false, "947÷2=473, 1", 2) | Out-Null
$d.Content.Find.Execute("107÷7=15, 2", $true, $false, $false, $false, $false, $true, 1, $false, "428÷6=71, 2", 2) | Out-Null
$d.Content.Find.Execute("868÷8=108, 4", $true, $false, $false, $false, $false, $true, 1, $false, "508÷6=84, 4", 2) | Out-Null
$d.Content.Find.Execute("447÷5=89, 2", $true, $false, $false, $false, $false, $true, 1, $false, "529÷6=88, 1", 2) | Out-Null
$d.Content.Find.Execute("153÷2=76, 1", $true, $false, $false, $false, $false, $true, 1, $false, "398÷5=79, 3", 2) | Out-Null
$d.Content.Find.Execute("139÷6=23, 1", $true, $false, $false, $false, $false, $true, 1, $false, "295÷8=36, 7", 2) | Out-Null
